$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room for two new item rows right after the existing row 9 (the
#    "بلاستر مترسيلك 2.5 سم" line). This pushes the totals row (old row 10)
#    down to row 12 and the footer row (old row 11) down to row 13.
# ---------------------------------------------------------------------------
$ws.Rows("10:11").Insert()

# Seed the two fresh rows by copying row 9 (still holding the original
# "بلاستر" data at this point) so they inherit the same cell styles and
# merged-cell layout as the other item rows.
$ws.Range("A9:Q9").Copy($ws.Range("A10:Q10"))
$ws.Range("A9:Q9").Copy($ws.Range("A11:Q11"))

# Restore/adjust row heights to match the new layout.
$ws.Rows("10").RowHeight = 24.75
$ws.Rows("11").RowHeight = 25.5
$ws.Rows("12").RowHeight = 25.5

# ---------------------------------------------------------------------------
# 2) Row 9 becomes the new "OTRIVIN" item (previously the "بلاستر" item).
#    Columns L and P carry numeric-looking text in this report, so flip the
#    cell to Text, write the string, then restore the original number
#    format so the style index is reused instead of minted anew.
# ---------------------------------------------------------------------------
$ws.Range("C9").Value = "OTRIVIN 0.1% ADULT NASAL DROPS 15 ML"
$ws.Range("H9").Value = "5:0"

$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "1"
$ws.Range("L9").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

$ws.Range("N9").Value = "24.00"

$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "24.0000"
$ws.Range("P9").NumberFormat = "0.00"
# Q9 ("1:0") is unchanged.

# ---------------------------------------------------------------------------
# 3) Row 10 becomes the new "PENTACOLD" item.
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "PENTACOLD SYRUP 120 ML"
$ws.Range("H10").Value = "0:0"

$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "1"
$ws.Range("L10").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

$ws.Range("N10").Value = "32.00"

$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value = "32.0000"
$ws.Range("P10").NumberFormat = "0.00"
# Q10 ("1:0") already correct from the row-9 copy.

# ---------------------------------------------------------------------------
# 4) Row 11 keeps the original "بلاستر مترسيلك 2.5 سم" data (copied from the
#    old row 9) - only the running item number changes.
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = 5

# ---------------------------------------------------------------------------
# 5) Update the grand-total cell (old row 10, now row 12): +24.00 +32.00
# ---------------------------------------------------------------------------
$ws.Range("P12").Value = 248.27000000000001

# ---------------------------------------------------------------------------
# 6) Footer row (old row 11, now row 13): bump the printed timestamp.
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "Sunday, 31 August, 2025 11:11 AM"
